# Coupling Config.xlsx - "Added look ahead and implemented a future co2 price"
#
# - Change 'End Year' (B3) from 2030 to 2025
# - Add a new parameter row: 'Look Ahead' (A4) = 7 (B4)
# - Add a threaded comment on B4 explaining the new parameter
# - Move the active selection to C6 (as left by the author after editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 'End Year' changes from 2030 to 2025
$ws.Range("B3").Value = 2025

# New row 4: 'Look Ahead' parameter
$ws.Range("A4").Value = "Look Ahead"
$ws.Range("B4").Value = 7

# Threaded comment on the new 'Look Ahead' value explaining its purpose
$excel.UserName = "Jim Hommes"
$comment = $ws.Range("B4").AddCommentThreaded("Be sure that there is data ready until 'End Year' + 'Look Ahead'")

# Leave the selection where the author left it after making the edit
$ws.Range("C6").Select()
